# Append the new simulation-result rows (124-159) produced by the smarter
# replacement algorithm, then move the viewport/selection down to where the
# new data was typed in (mirrors what a person scrolling to row ~121 and
# clicking I123 after pasting the new block would produce).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Tab-separated-looking, but we use commas: A,B,C,D,E,F,G,H for rows 124..159.
# Row 123 is intentionally left blank - the source data has a gap there.
$newData = @"
0.5,0,4,5.7875370370370305E-4,4,0.29971483988425801,0.19492424740740599,0.17097107849536899
0.5,4,4,5.7875370370370305E-4,4,0.57823282537037202,0.33600269710648001,0.251273154884257
0.5,8,4,5.7875370370370305E-4,4,0.97179981273148996,0.52105919386573696,0.38930591321758901
0.5,16,4,5.7875370370370305E-4,4,2.0033342421065399,1.02873470833335,0.66243425483796803
0.5,0,16,5.7875370370370305E-4,4,7.8507939907407107E-2,6.55944978935184E-2,5.9843132962962897E-2
0.5,16,16,5.7875370370370305E-4,4,1.04169879129634,0.31487818692129299,0.25041225874999701
0.5,32,16,5.7875370370370305E-4,4,1.4602245321297,0.66974825476854505,0.44397643495370998
0.5,64,16,5.7875370370370305E-4,4,6.7798970750898304,1.3549492334258599,0.74276526891204098
0.5,0,64,5.7875370370370305E-4,4,2.7599317245370299E-2,2.4191904814814799E-2,1.76881600694444E-2
0.5,64,64,5.7875370370370305E-4,4,1.03612105247665,0.302398810185182,0.225149659583331
0.5,128,64,5.7875370370370305E-4,4,5.0064365696783302,2.5155529731466801,0.45670178201388301
0.5,256,64,5.7875370370370305E-4,4,73.214955144815804,8.52171422177579,1.5470520565290899
0.5,0,4,5.7875370370370305E-4,16,0.25082462076388801,0.202238247337962,0.168815220949073
0.5,4,4,5.7875370370370305E-4,16,0.61168478944444404,0.33897604425925598,0.24599926175925699
0.5,8,4,5.7875370370370305E-4,16,1.2781775546296399,0.53110780504629096,0.379980744166663
0.5,16,4,5.7875370370370305E-4,16,4.3294177215051901,1.10253303997689,0.75965764263889801
0.5,0,16,5.7875370370370305E-4,16,8.9236586689814601E-2,8.2711138680555393E-2,5.1089483194444403E-2
0.5,16,16,5.7875370370370305E-4,16,1.70628166925944,0.377332945972226,0.27754857303240998
0.5,32,16,5.7875370370370305E-4,16,8.3326353867575396,0.48997288555555102,0.52439426208335904
0.5,64,16,5.7875370370370305E-4,16,23.596685568222401,6.48037033015859,1.9062844804162999
0.5,0,64,5.7875370370370305E-4,16,2.3757839537037001E-2,1.7731566597222199E-2,1.78111452314814E-2
0.5,64,64,5.7875370370370305E-4,16,6.5426804007908999,0.29914332060184901,0.225945445925924
0.5,128,64,5.7875370370370305E-4,16,98.9561329835891,1.62032950984103,1.35553522155338
0.5,256,64,5.7875370370370305E-4,16,261.92280178671399,213.649078346192,19.2613066836625
0.5,0,4,5.7875370370370305E-4,64,0.27701322585648103,0.20126883488425701,0.174631695671295
0.5,4,4,5.7875370370370305E-4,64,0.752271298495366,0.384654180324071,0.27900269171296099
0.5,8,4,5.7875370370370305E-4,64,3.0110818942594499,2.1010495705555399,0.46214930124999398
0.5,16,4,5.7875370370370305E-4,64,6.06406555666765,4.5111391500468896,3.82927323960689
0.5,0,16,5.7875370370370305E-4,64,8.5084028865740499E-2,7.87539102314814E-2,4.7016504004629603E-2
0.5,16,16,5.7875370370370305E-4,64,9.3659350149287999,2.1138255585649302,0.844539107708416
0.5,32,16,5.7875370370370305E-4,64,23.907027772990698,20.029920539778001,17.291851236401801
0.5,64,16,5.7875370370370305E-4,64,69.997641602610102,57.206707278890299,54.261747995264898
0.5,0,64,5.7875370370370305E-4,64,2.3085038356481399E-2,2.1884124421296201E-2,2.0719382592592501E-2
0.5,64,64,5.7875370370370305E-4,64,107.93493241173,8.8984467107919194,1.7038364348659301
0.5,128,64,5.7875370370370305E-4,64,219.544538480131,231.16270800188599,227.643465886934
0.5,256,64,5.7875370370370305E-4,64,477.00387594336701,459.159025543523,490.76132895169297
"@

$lines = $newData.Split("`n")
$startRow = 124
$r = $startRow
foreach ($line in $lines) {
    $t = $line.Trim()
    if ($t.Length -eq 0) { continue }
    $fields = $t.Split(",")
    for ($c = 0; $c -lt $fields.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = [double]$fields[$c]
    }
    $r = $r + 1
}
$endRow = $r - 1

# Column D holds the (tiny, constant) timestep and is formatted in
# scientific notation in the existing data - carry that formatting onto the
# freshly written cells too.
$ws.Range("D" + $startRow + ":D" + $endRow).NumberFormat = "0.00E+00"

# Scroll the view down to the newly added block and leave the selection
# where it would land right after typing/pasting the new rows.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 121
$ws.Range("I123").Select()
